$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.281.38"
$ws.Range("E2").Value = "  -1.01%  "

$ws.Range("D3").Value = "3.243.49"
$ws.Range("E3").Value = "  +3.15%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "594.68"
$ws.Range("E5").Value = "  -1.05%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "140.37"
$ws.Range("E6").Value = "  -1.10%  "

$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("D8").Value = "3.238.62"
$ws.Range("E8").Value = "  +3.11%  "

$ws.Range("E9").Value = "  -2.10%  "

$ws.Range("E10").Value = "  -1.06%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.35"
$ws.Range("E11").Value = "  -0.48%  "

$ws.Range("E12").Value = "  -0.45%  "

$ws.Range("E13").Value = "  -2.74%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.24"
$ws.Range("E14").Value = "  -1.81%  "

$ws.Range("D15").Value = "3.770.76"
$ws.Range("E15").Value = "  +2.86%  "

$ws.Range("E16").Value = "  -0.22%  "

$ws.Range("D17").Value = "3.239.16"
$ws.Range("E17").Value = "  +3.13%  "

$ws.Range("D18").Value = "63.322.95"
$ws.Range("E18").Value = "  -0.97%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.75"
$ws.Range("E19").Value = "  -1.00%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "472.85"
$ws.Range("E20").Value = "  -2.96%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.14"

$ws.Range("E22").Value = "  +2.66%  "

$ws.Range("E23").Value = "  +2.67%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.71"
$ws.Range("E24").Value = "  -5.17%  "

$ws.Range("E25").Value = "  -0.47%  "

$ws.Range("E26").Value = "  -0.05%  "

$ws.Range("E27").Value = "  -1.16%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.37"
$ws.Range("E28").Value = "  +5.77%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.09"
$ws.Range("E29").Value = "  -1.20%  "

$ws.Range("E30").Value = "  +2.91%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "27.42"
$ws.Range("E31").Value = "  -0.07%  "

$ws.Range("E32").Value = "  -0.14%  "

$ws.Range("E33").Value = "  -4.06%  "

$ws.Range("E34").Value = "  -4.44%  "

$ws.Range("E35").Value = "  -1.02%  "

$ws.Range("E36").Value = "  -2.19%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "52.65"
$ws.Range("E37").Value = "  +0.00%  "

$ws.Range("D38").Value = "0.0₃0710"
$ws.Range("E38").Value = "  -4.67%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0392"
$ws.Range("E39").Value = "  -1.11%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "422.47"
$ws.Range("E40").Value = "  -2.03%  "

$ws.Range("E41").Value = "  -0.01%  "

$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "2.970.80"
$ws.Range("E42").Value = "  +1.98%  "

$ws.Range("B43").Value = "dogwifhat"
$ws.Range("C43").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.74"
$ws.Range("E43").Value = "  -6.18%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.110"
$ws.Range("E44").Value = "  -8.76%  "

$ws.Range("E45").Value = "  +2.64%  "

$ws.Range("E46").Value = "  -0.92%  "

$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "25.88"
$ws.Range("E48").Value = "  +0.64%  "

$ws.Range("B49").Value = "ThetaToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.33"
$ws.Range("E49").Value = "  -2.70%  "

$ws.Range("E50").Value = "  -0.38%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "121.54"
$ws.Range("E51").Value = "  +0.65%  "
